# #5: property building done
# Append a second "property" row to both the 汽車 (vehicle) and 債務 (debt)
# sheets, duplicating the first data row's values but tagging each new
# row with a leading sequence number in column A.

$wb = $excel.ActiveWorkbook

# --- Sheet 1 (汽車) ---------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A2").Value = 29
$ws1.Range("B2").Value = "中華FB308W"
$ws1.Range("C2").Value = 2835
$ws1.Range("D2").Value = "陳亭妃"
$ws1.Range("E2").Value = "98年03月24日"
$ws1.Range("F2").Value = "繼承"
$ws1.Range("G2").Value = 10000

# Stamp the new row with its own (blank/default) style so it gets its
# own cellXfs entry, distinct from the header row's bordered style.
$ws1.Range("A2:G2").HorizontalAlignment = 1

# --- Sheet 2 (債務) ---------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A2").Value = 83
$ws2.Range("B2").Value = "中期放款"
$ws2.Range("C2").Value = "陳亭妃"
$ws2.Range("D2").Value = "合作金庫商業銀行臺南市北區曲門路"
$ws2.Range("E2").Value = 930000
$ws2.Range("F2").Value = "89年03月29日"
$ws2.Range("G2").Value = "信用貸款"

$ws2.Range("A2:G2").HorizontalAlignment = 1
